# Error Calculations and Plots
# Apply corrected/updated values across the data table and remove the
# last two rows (the data for "RM 232" is dropped and subsequent rows
# shift up by one; rows that mapped to SC193/SC232 duplicates at the
# bottom are removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two trailing rows first so everything below row 26 shifts
# up by one (RM 232 row removed; SC5..SC232 move up one row).
$ws.Rows.Item(35).Delete()
$ws.Rows.Item(34).Delete()

# Now write the corrected values for rows 2-33.
$ws.Range("A2").Value = "RM 2"
$ws.Range("B2").Value = -19.7
$ws.Range("C2").Value = 14.9
$ws.Range("D2").Value = -13.5
$ws.Range("E2").Value = -7.2
$ws.Range("F2").Value = 18.03

$ws.Range("A3").Value = "RM 8"
$ws.Range("B3").Value = -19.7
$ws.Range("C3").Value = 11.2
$ws.Range("D3").Value = -14.2
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()

$ws.Range("A4").Value = "RM 9"
$ws.Range("B4").Value = -18.7
$ws.Range("C4").Value = 11
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = -6.4
$ws.Range("F4").Value = 17.97

$ws.Range("A5").Value = "RM 14"
$ws.Range("B5").Value = -19.5
$ws.Range("C5").Value = 12.3
$ws.Range("D5").Value = -14.4
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 17.66

$ws.Range("A6").Value = "RM 21"
$ws.Range("B6").Value = -19.8
$ws.Range("C6").Value = 15.1
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = -5.7
$ws.Range("F6").Value = 16.43

$ws.Range("A7").Value = "RM 32"
$ws.Range("B7").Value = -19.5
$ws.Range("C7").Value = 15
$ws.Range("D7").Value = -13.8
$ws.Range("E7").Value = -7.1
$ws.Range("F7").Value = 17.24

$ws.Range("A8").Value = "RM 38"
$ws.Range("B8").Value = -19.9
$ws.Range("C8").Value = 15.5
$ws.Range("D8").Value = -13.9
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()

$ws.Range("A9").Value = "RM 42"
$ws.Range("B9").Value = -20.6
$ws.Range("C9").Value = 10.5
$ws.Range("D9").Value = -14.5
$ws.Range("E9").ClearContents()
$ws.Range("F9").ClearContents()

$ws.Range("A10").Value = "RM 52 a"
$ws.Range("B10").Value = -19.8
$ws.Range("C10").Value = 11.5
$ws.Range("D10").Value = -14.7
$ws.Range("E10").Value = -6.1
$ws.Range("F10").Value = 16.43

$ws.Range("A11").Value = "RM 58"
$ws.Range("B11").Value = -20.8
$ws.Range("C11").ClearContents()
$ws.Range("D11").Value = -15.5
$ws.Range("E11").Value = -7.9
$ws.Range("F11").Value = 17.65

$ws.Range("A12").Value = "RM 81"
$ws.Range("B12").Value = -19.9
$ws.Range("C12").ClearContents()
$ws.Range("D12").Value = -14.1
$ws.Range("E12").Value = -5.3
$ws.Range("F12").Value = 17.45

$ws.Range("A13").Value = "RM 88"
$ws.Range("B13").Value = -19.9
$ws.Range("C13").Value = 12.5
$ws.Range("D13").Value = -13.9
$ws.Range("E13").Value = -5.3
$ws.Range("F13").Value = 17.1

$ws.Range("A14").Value = "RM 90"
$ws.Range("B14").Value = -19.6
$ws.Range("C14").Value = 14.4
$ws.Range("D14").ClearContents()
$ws.Range("E14").Value = -5.4
$ws.Range("F14").Value = 17.76

$ws.Range("A15").Value = "RM 95"
$ws.Range("B15").Value = -19.1
$ws.Range("C15").ClearContents()
$ws.Range("D15").Value = -15.2
$ws.Range("E15").Value = -8.4
$ws.Range("F15").Value = 16.2

$ws.Range("A16").Value = "RM 103"
$ws.Range("B16").Value = -19.5
$ws.Range("C16").Value = 13.5
$ws.Range("D16").Value = -15.3
$ws.Range("E16").Value = -6.9
$ws.Range("F16").Value = 17.34

$ws.Range("A17").Value = "RM 116"
$ws.Range("B17").Value = -19.4
$ws.Range("C17").Value = 11.2
$ws.Range("D17").Value = -14.7
$ws.Range("E17").Value = -7.3
$ws.Range("F17").Value = 17.78

$ws.Range("A18").Value = "RM 120"
$ws.Range("B18").Value = -19.6
$ws.Range("C18").Value = 11.5
$ws.Range("D18").Value = -15.2
$ws.Range("E18").Value = -8.5
$ws.Range("F18").Value = 18.35

$ws.Range("A19").Value = "RM 125"
$ws.Range("B19").Value = -20.6
$ws.Range("C19").Value = 13.2
$ws.Range("D19").Value = -15.5
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()

$ws.Range("A20").Value = "RM 134"
$ws.Range("B20").Value = -19
$ws.Range("C20").Value = 12.5
$ws.Range("D20").Value = -14
$ws.Range("E20").Value = -7.2
$ws.Range("F20").Value = 17.73

$ws.Range("A21").Value = "RM 135"
$ws.Range("B21").Value = -18.9
$ws.Range("C21").Value = 12.7
$ws.Range("D21").Value = -14.3
$ws.Range("E21").Value = -8.7
$ws.Range("F21").Value = 16.58

$ws.Range("A22").Value = "RM 138"
$ws.Range("B22").Value = -19.3
$ws.Range("C22").Value = 13.5
$ws.Range("D22").ClearContents()
$ws.Range("E22").ClearContents()
$ws.Range("F22").Value = 16.81

$ws.Range("A23").Value = "RM 140"
$ws.Range("B23").Value = -19.5
$ws.Range("C23").Value = 12.2
$ws.Range("D23").ClearContents()
$ws.Range("E23").Value = -7
$ws.Range("F23").Value = 16.48

$ws.Range("A24").Value = "RM 142a"
$ws.Range("B24").Value = -17.7
$ws.Range("C24").Value = 12.7
$ws.Range("D24").ClearContents()
$ws.Range("E24").Value = -8.1
$ws.Range("F24").Value = 16.78

$ws.Range("A25").Value = "RM 145"
$ws.Range("B25").Value = -19.5
$ws.Range("C25").Value = 10.7
$ws.Range("D25").Value = -15.5
$ws.Range("E25").Value = -7.1
$ws.Range("F25").Value = 16.6

$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").Value = -20.2
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38

$ws.Range("A27").Value = "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("C27").ClearContents()
$ws.Range("D27").Value = -14.6
$ws.Range("E27").ClearContents()
$ws.Range("F27").Value = 17

$ws.Range("A28").Value = "SC 105"
$ws.Range("B28").Value = -19.6
$ws.Range("C28").ClearContents()
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = 17.44

$ws.Range("A29").Value = "SC 119"
$ws.Range("B29").Value = -19.5
$ws.Range("C29").Value = 11.2
$ws.Range("D29").Value = -13
$ws.Range("E29").ClearContents()
$ws.Range("F29").Value = 18.06

$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89

$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").ClearContents()
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").ClearContents()

$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").Value = -19.9
$ws.Range("C32").ClearContents()
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53
